$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("D4").Value = 21.3
$ws.Range("B11").Value = 0.6
$ws.Range("G11").Value = "HolaCambiado2"

# Add the new defined names (workbook scope)
$wb.Names.Add("PUE.SLIDE.0.1.0.1.porcentaje1", "=Sheet1!`$B`$11")
$wb.Names.Add("PUE.STRING.cadena1", "=Sheet1!`$G`$11")

# Update the selected cell/range to match the saved view state
$ws.Range("B11").Select()
